$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of the
# existing header cells (bold, bordered, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new data columns I and J for rows 2-18.
$data = @(
    @(7, 8),
    @(1, 3),
    @(7, 8),
    @(1, 5),
    @(4, 7),
    @(8, 8),
    @(5, 7),
    @(7, 8),
    @(6, 7),
    @(8, 8),
    @(2, 5),
    @(5, 6),
    @(8, 9),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(6, 7)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 9).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $r = $r + 1
}
